# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.606.38"
$ws.Range("E2").Value = "'  +1.21%  "
$ws.Range("D3").Value = "'1.801.68"
$ws.Range("E3").Value = "'  +1.03%  "
$ws.Range("E4").Value = "'  -0.26%  "
$ws.Range("D5").Value = "'227.57"
$ws.Range("E5").Value = "'  +0.63%  "
$ws.Range("E6").Value = "'  +1.77%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.23%  "
$ws.Range("D8").Value = "'32.74"
$ws.Range("E8").Value = "'  +3.06%  "
$ws.Range("D9").Value = "'0.297"
$ws.Range("E9").Value = "'  +1.73%  "
$ws.Range("E10").Value = "'  +0.89%  "
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = "'  +0.32%  "
$ws.Range("D12").Value = "'2.061.09"
$ws.Range("E12").Value = "'  +0.95%  "
$ws.Range("D13").Value = "'11.18"
$ws.Range("E13").Value = "'  +2.19%  "
$ws.Range("D14").Value = "'1.794.02"
$ws.Range("E14").Value = "'  +0.76%  "
$ws.Range("D15").Value = "'0.641"
$ws.Range("E15").Value = "'  +2.80%  "
$ws.Range("D16").Value = "'34.573.04"
$ws.Range("E16").Value = "'  +1.18%  "
$ws.Range("D17").Value = "'4.34"
$ws.Range("E17").Value = "'  +3.72%  "
$ws.Range("D18").Value = "'68.91"
$ws.Range("E18").Value = "'  +1.48%  "
$ws.Range("D19").Value = "'0.0₃0806"
$ws.Range("E19").Value = "'  +0.79%  "
$ws.Range("D20").Value = "'247.60"
$ws.Range("E20").Value = "'  +0.54%  "
$ws.Range("D21").Value = "'11.34"
$ws.Range("E21").Value = "'  +3.42%  "
$ws.Range("E22").Value = "'  -0.13%  "
$ws.Range("E23").Value = "'  +2.75%  "
$ws.Range("D24").Value = "'169.66"
$ws.Range("E24").Value = "'  +4.57%  "
$ws.Range("E25").Value = "'  +1.60%  "
$ws.Range("E26").Value = "'  +1.53%  "
$ws.Range("D27").Value = "'16.61"
$ws.Range("E27").Value = "'  +1.98%  "
$ws.Range("E28").Value = "'  +2.30%  "
$ws.Range("E29").Value = "'  -0.43%  "
$ws.Range("D30").Value = "'4.14"
$ws.Range("E30").Value = "'  +11.08%  "
$ws.Range("B31").Value = "'Hedera"
$ws.Range("C31").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0527"
$ws.Range("E31").Value = "'  +1.34%  "
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.83"
$ws.Range("E32").Value = "'  +2.59%  "
$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "'  +0.75%  "
$ws.Range("E34").Value = "'  +3.09%  "
$ws.Range("D35").Value = "'1.432.46"
$ws.Range("E35").Value = "'  -0.77%  "
$ws.Range("D36").Value = "'2.62"
$ws.Range("E36").Value = "'  +9.17%  "
$ws.Range("D37").Value = "'0.679"
$ws.Range("E37").Value = "'  +3.81%  "
$ws.Range("E38").Value = "'  +3.29%  "
$ws.Range("E39").Value = "'  +0.52%  "
$ws.Range("D40").Value = "'85.39"
$ws.Range("E40").Value = "'  +6.44%  "
$ws.Range("E41").Value = "'  +2.82%  "
$ws.Range("E42").Value = "'  +2.17%  "
$ws.Range("E43").Value = "'  +3.28%  "
$ws.Range("D44").Value = "'13.86"
$ws.Range("E44").Value = "'  +2.77%  "
$ws.Range("E45").Value = "'  +3.00%  "
$ws.Range("D46").Value = "'6.11"
$ws.Range("E46").Value = "'  +0.52%  "
$ws.Range("E47").Value = "'  +0.72%  "
$ws.Range("D48").Value = "'1.960.25"
$ws.Range("E48").Value = "'  +0.86%  "
$ws.Range("D49").Value = "'105.92"
$ws.Range("E49").Value = "'  +1.44%  "
$ws.Range("E50").Value = "'  -0.23%  "
$ws.Range("E51").Value = "'  -4.68%  "
